$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows above row 98; this pushes the existing rows 98-108
# down to 100-110 and extends the used range to A1:T110.
$ws.Rows.Item(98).Resize(2).Insert()

# Row 98: new weekly record (Kiwi, Primera, Provincia de Curico)
$ws.Cells.Item(98, 1).Value = 11
$ws.Cells.Item(98, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(98, 3).Value = "Bíobío"
$ws.Cells.Item(98, 4).Value = 44491
$ws.Cells.Item(98, 5).Value = 8
$ws.Cells.Item(98, 6).Value = "Fruta"
$ws.Cells.Item(98, 7).Value = 100101
$ws.Cells.Item(98, 8).Value = "Berries"
$ws.Cells.Item(98, 9).Value = 100101007
$ws.Cells.Item(98, 10).Value = "Kiwi"
$ws.Cells.Item(98, 11).Value = "Hayward"
$ws.Cells.Item(98, 12).Value = "Primera"
$ws.Cells.Item(98, 13).Value = 80
$ws.Cells.Item(98, 14).Value = 17000
$ws.Cells.Item(98, 15).Value = 17000
$ws.Cells.Item(98, 16).Value = 17000
$ws.Cells.Item(98, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(98, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(98, 19).Value = 944
$ws.Cells.Item(98, 20).Value = 18

# Row 99: new weekly record (Kiwi, Segunda, Provincia de Curico)
$ws.Cells.Item(99, 1).Value = 11
$ws.Cells.Item(99, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(99, 3).Value = "Bíobío"
$ws.Cells.Item(99, 4).Value = 44491
$ws.Cells.Item(99, 5).Value = 8
$ws.Cells.Item(99, 6).Value = "Fruta"
$ws.Cells.Item(99, 7).Value = 100101
$ws.Cells.Item(99, 8).Value = "Berries"
$ws.Cells.Item(99, 9).Value = 100101007
$ws.Cells.Item(99, 10).Value = "Kiwi"
$ws.Cells.Item(99, 11).Value = "Hayward"
$ws.Cells.Item(99, 12).Value = "Segunda"
$ws.Cells.Item(99, 13).Value = 170
$ws.Cells.Item(99, 14).Value = 14000
$ws.Cells.Item(99, 15).Value = 15000
$ws.Cells.Item(99, 16).Value = 14412
$ws.Cells.Item(99, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(99, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(99, 19).Value = 801
$ws.Cells.Item(99, 20).Value = 18
